$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "slope" column (column L) entirely; cells to the right shift left.
$ws.Columns("L").Delete()

# Rename the "delta_5" header (now still in column I) to "delta".
$ws.Range("I1").Value = "delta"
